$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D, E and G columns keep their text representation (the source data
# are formatted strings, not numbers) by pre-setting a text number format,
# writing the values, then reverting to the Normal style so no stray
# per-cell style index is left on the written cells.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "301.32"
$ws.Range("E2").Value = "-0.59%"
$ws.Range("G2").Value = "19"
$ws.Range("D3").Value = "31.39"
$ws.Range("E3").Value = "-1.68%"
$ws.Range("G3").Value = "19"
$ws.Range("D4").Value = "5.156"
$ws.Range("E4").Value = "-1.87%"
$ws.Range("G4").Value = "19"
$ws.Range("D5").Value = "0.07370"
$ws.Range("E5").Value = "-1.17%"
$ws.Range("G5").Value = "19"
$ws.Range("D6").Value = "2.428"
$ws.Range("E6").Value = "58.71%"
$ws.Range("G6").Value = "19"
$ws.Range("D7").Value = "7.924"
$ws.Range("E7").Value = "0.95%"
$ws.Range("G7").Value = "19"
$ws.Range("E8").Value = "-0.99%"
$ws.Range("G8").Value = "19"
$ws.Range("D9").Value = "0.9215"
$ws.Range("E9").Value = "0.32%"
$ws.Range("G9").Value = "19"
$ws.Range("D10").Value = "0.1740"
$ws.Range("E10").Value = "3.60%"
$ws.Range("G10").Value = "19"
$ws.Range("D11").Value = "0.07462"
$ws.Range("E11").Value = "-7.22%"
$ws.Range("G11").Value = "19"
$ws.Range("D12").Value = "0.08135"
$ws.Range("E12").Value = "2.31%"
$ws.Range("G12").Value = "19"
$ws.Range("D13").Value = "0.03045"
$ws.Range("E13").Value = "0.40%"
$ws.Range("G13").Value = "19"
$ws.Range("D14").Value = "0.09943"
$ws.Range("E14").Value = "0.50%"
$ws.Range("G14").Value = "19"
$ws.Range("D15").Value = "0.001498"
$ws.Range("E15").Value = "-0.03%"
$ws.Range("G15").Value = "19"
$ws.Range("D16").Value = "0.006097"
$ws.Range("E16").Value = "-5.06%"
$ws.Range("G16").Value = "19"
$ws.Range("D17").Value = "3.453"
$ws.Range("E17").Value = "-0.25%"
$ws.Range("G17").Value = "19"
$ws.Range("D18").Value = "2.232"
$ws.Range("E18").Value = "0.09%"
$ws.Range("G18").Value = "19"
$ws.Range("D19").Value = "0.3290"
$ws.Range("E19").Value = "-1.05%"
$ws.Range("G19").Value = "19"
$ws.Range("D20").Value = "0.1338"
$ws.Range("E20").Value = "0.00%"
$ws.Range("G20").Value = "19"
$ws.Range("D21").Value = "4.659"
$ws.Range("E21").Value = "4.08%"
$ws.Range("G21").Value = "19"
$ws.Range("D22").Value = "0.04635"
$ws.Range("E22").Value = "0.80%"
$ws.Range("G22").Value = "19"
$ws.Range("G23").Value = "19"
$ws.Range("D24").Value = "0.001224"
$ws.Range("E24").Value = "0.72%"
$ws.Range("G24").Value = "19"
$ws.Range("E25").Value = "0.78%"
$ws.Range("G25").Value = "19"
$ws.Range("E26").Value = "-7.03%"
$ws.Range("G26").Value = "19"
$ws.Range("E27").Value = "7.09%"
$ws.Range("G27").Value = "19"
$ws.Range("G28").Value = "19"
$ws.Range("G29").Value = "19"
$ws.Range("G30").Value = "19"
$ws.Range("G31").Value = "19"
$ws.Range("G32").Value = "19"
$ws.Range("G33").Value = "19"
$ws.Range("G34").Value = "19"
$ws.Range("G35").Value = "19"
$ws.Range("G36").Value = "19"
$ws.Range("G37").Value = "19"
$ws.Range("G38").Value = "19"
$ws.Range("D39").Value = "0.01733"
$ws.Range("E39").Value = "-0.18%"
$ws.Range("G39").Value = "19"
$ws.Range("E40").Value = "0.86%"
$ws.Range("G40").Value = "19"
$ws.Range("D41").Value = "0.007256"
$ws.Range("E41").Value = "1.07%"
$ws.Range("G41").Value = "19"
$ws.Range("D42").Value = "0.1348"
$ws.Range("E42").Value = "0.10%"
$ws.Range("G42").Value = "19"
$ws.Range("D43").Value = "0.002222"
$ws.Range("E43").Value = "1.04%"
$ws.Range("G43").Value = "19"
$ws.Range("D44").Value = "0.01089"
$ws.Range("E44").Value = "-15.14%"
$ws.Range("G44").Value = "19"
$ws.Range("D45").Value = "0.00006289"
$ws.Range("E45").Value = "2.30%"
$ws.Range("G45").Value = "19"
$ws.Range("E46").Value = "-22.98%"
$ws.Range("G46").Value = "19"
$ws.Range("D47").Value = "1.928"
$ws.Range("E47").Value = "171.71%"
$ws.Range("G47").Value = "19"
$ws.Range("G48").Value = "19"
$ws.Range("G49").Value = "19"
$ws.Range("G50").Value = "19"
$ws.Range("G51").Value = "19"

# Revert number format to the workbook default style so the text cells
# keep rendering exactly like their neighbours.
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("E2:E51").Style = "Normal"
$ws.Range("G2:G51").Style = "Normal"
